# RTI 7 Scenarios input files checkin for employee changes
# Update the "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 6" marker text
# to "... EMPLOYEE 106" across all sheets that reference it.

$wb = $excel.ActiveWorkbook

$oldText = "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 6"
$newText = "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 106"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq $oldText) {
            $cell.Value2 = $newText
        }
    }
}
